$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new headers for columns E and F
$ws.Range("E1").Value = "start"
$ws.Range("F1").Value = "stop"

# Fill columns E and F with value 20 for all data rows (row 2 through the last used row)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = 20
    $ws.Cells.Item($r, 6).Value = 20
}

# Update the view state to match the author's final cursor position/selection
$ws.Range("A15").Select()
$ws.Range("G66").Select()
